# Aggiornamento fino a 28/06 incluso
# Appends 32 new daily rows (270-301) after the last existing row (269),
# continuing the date series (serial 44344 .. 44375) with 0 values in
# columns B, C, D - matching the formatting (style) of the preceding row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 269
$firstNewRow = 270
$lastNewRow = 301

# Copy the formatting of the last existing data row down onto the new
# block so the new date cells (column A) pick up the same style (border,
# bold, centered, date number format) as the existing ones, without
# creating any new style entries.
$ws.Range("A$lastRow`:D$lastRow").Copy() | Out-Null
$ws.Range("A$firstNewRow`:D$lastNewRow").PasteSpecial(-4122) | Out-Null

$startSerial = 44344
for ($r = $firstNewRow; $r -le $lastNewRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $startSerial + ($r - $firstNewRow)
    $ws.Cells.Item($r, 2).Value = 0
    $ws.Cells.Item($r, 3).Value = 0
    $ws.Cells.Item($r, 4).Value = 0
}
